# Updates to the contact list on the "US Groups" sheet:
# - new "contacted on" / "by" / "response received" header cells (F1:H1, bold)
# - contact dates + initials for the rows that got a response
# - a note about the new outline .tex document (KPiX_RD_Status.docx) on row 17
# Also restores the view state (active sheet/selection) for every sheet to
# match what was left after the edit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("US Groups")

# --- new header cells in row 1 --------------------------------------------
$ws.Range("F1").Value = "contacted on"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").WrapText = $false

$ws.Range("G1").Value = "by"
$ws.Range("G1").Font.Bold = $true

$ws.Range("H1").Value = "response received"
$ws.Range("H1").Font.Bold = $true

# --- Alignment group (row 2-3 block) contacted 13-Mar by Jan --------------
$contacted1 = Get-Date -Year 2014 -Month 3 -Day 13 -Hour 0 -Minute 0 -Second 0
$ws.Range("F3").Value = $contacted1
$ws.Range("F3").NumberFormat = "d-mmm"
$ws.Range("G3").Value = "Jan"

# --- GEM calorimetry block contacted 13-Mar by Jan -------------------------
$ws.Range("F7").Value = $contacted1
$ws.Range("F7").NumberFormat = "d-mmm"
$ws.Range("G7").Value = "Jan"

# --- Kpix row: contacted 13-Mar by Jan, response received 18-Mar ----------
$ws.Range("F17").Value = $contacted1
$ws.Range("F17").NumberFormat = "d-mmm"
$ws.Range("G17").Value = "Jan"

$response1 = Get-Date -Year 2014 -Month 3 -Day 18 -Hour 0 -Minute 0 -Second 0
$ws.Range("H17").Value = $response1
$ws.Range("H17").NumberFormat = "d-mmm"
$ws.Range("I17").Value = "KPiX_RD_Status.docx"

# --- restore selections on the other sheets (view-only state) -------------
$wsTPC = $wb.Worksheets.Item("TPC")
$wsTPC.Range("E24").Select()

$wsVertex = $wb.Worksheets.Item("Vertex Detector")
$wsVertex.Range("A4").Select()

$wsCalice = $wb.Worksheets.Item("Calice")
$wsCalice.Range("B8").Select()

# "US Groups" becomes the active tab; leave its selection on H17, the last
# cell touched by the update.
$ws.Activate()
$ws.Range("H17").Select()
